# MM_26 Bill of Parts.xlsx - add Wheels, IR Sensor Module, and Dual H-Bridge parts
# (kicad / voltage regulator additions per commit message) as rows 4-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4: Wheels
# ---------------------------------------------------------------------------
# Start from row 3's formatting (same column layout/styles as other part rows)
$ws.Range("A3:I3").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122)
$ws.Range("C4").Clear()
# I column here is plain text (no hyperlink relationship), matching I2's style
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Wheels"
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 8.9499999999999993
$ws.Range("E4").Value = 1090
$ws.Range("F4").Formula = "=HYPERLINK(I4)"
$ws.Range("G4").Value = "attach to motor bracket"
$ws.Range("H4").Formula = "=(B4+C4)*D4"
$ws.Range("I4").Value = "https://www.pololu.com/product/1090"

# ---------------------------------------------------------------------------
# Row 5: IR Sensor Module (sparse row - just part name and quantity)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "IR Sensor Module"
$ws.Range("B5").Value = 5

# ---------------------------------------------------------------------------
# Row 6: Dual H-Bridge
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Dual H-Bridge"
$ws.Range("B6").Value = 1
$ws.Range("D6").Value = 7.69
$ws.Range("E6").Value = "497-1390-5-ND"
$ws.Range("G6").Value = "IC dual-H bridge"
$ws.Range("I6").Value = "https://www.digikey.com/en/products/detail/stmicroelectronics/L293DD/585913"
# This one keeps a real hyperlink relationship (like row 3's I3)
[void]$ws.Hyperlinks.Add($ws.Range("I6"), "https://www.digikey.com/en/products/detail/stmicroelectronics/L293DD/585913")

# Apply the row-3 style template to row 6 as well, then clean up C6 / reapply formulas
$ws.Range("A3:I3").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$ws.Range("C6").Clear()

$ws.Range("F6").Formula = "=HYPERLINK(I6)"
$ws.Range("H6").Formula = "=(B6+C6)*D6"

# ---------------------------------------------------------------------------
# Final selection, matching the saved-file cursor position
# ---------------------------------------------------------------------------
[void]$ws.Range("G11").Select()
